$p = $ppt.ActivePresentation

# --- 1) Switch the three tables (slides 14-16) to the "Medium Style 2 -
#        Accent 1" built-in table style. ---------------------------------
$newTableStyleId = "{340B9531-C0DE-4655-9DFE-57082375668E}"
foreach ($slideIdx in 14..16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($shpIdx = 1; $shpIdx -le $slide.Shapes.Count; $shpIdx++) {
        $shape = $slide.Shapes.Item($shpIdx)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId, $false)
        }
    }
}

# --- 2) Re-theme the deck from the "Integral" (Red Violet) design back
#        to the stock "Office Theme" colour scheme. ----------------------
$officeThemeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeThemeColors[$i - 1]
}

Write-Host "Applied table style + theme colour updates"
